# Release mCSD 3.9.0 with CP integrated
# Update the "Metadata" sheet of the CodeSystem-MCSDEndpointIdentifierTypes
# workbook to reflect the new version / date / contact / jurisdiction
# metadata for the 3.9.0 release.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 3.8.0 -> 3.9.0
$ws.Range("B3").Value = "3.9.0"

# Experimental: (empty) -> false
# A bare "false" is auto-typed as a Boolean by Excel's COM layer; prefixing
# with an apostrophe forces it to be entered (and stored) as literal text,
# matching the rest of this text-only metadata sheet. The leading
# apostrophe itself is not stored as part of the cell's text.
$ws.Range("B7").Value = "'false"
# Restore the plain "top/wrap" cell style (the quote-prefix entry above
# nudges Excel into minting a distinct, but equivalent, style record) by
# copying the formatting from a neighboring metadata-value cell.
$ws.Range("B8").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date: 2022-08-12T09:44:57-05:00 -> 2024-12-02T17:05:26-06:00
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact rows (10-12) now carry three distinct contact detail strings
# instead of the single repeated "No display for ContactDetail" text.
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"
